$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    "20ФиПЛ-1",
    "20ФиПЛ-2",
    "20ФИЛ-1",
    "20ФИЛ-2",
    "23ФИЛ-1",
    "23ФИЛ-2",
    "23ФиПЛ-1",
    "23ФиПЛ-2",
    "22ФИЛ-1",
    "22ФИЛ-2",
    "22ФиПЛ-1",
    "22ФиПЛ-2",
    "21ФИЛ-1",
    "21ФИЛ-2",
    "21ФиПЛ-1",
    "21ФиПЛ-2"
)

$baseUrl = "https://docs.google.com/spreadsheets/d/1E80xi4hgdd2JCox_hBIUjM08F5G6lRTfCC3E5K_7RoA/edit"
$gid = "gid=615601934"
$displayUrl = "$baseUrl#$gid"

$r = 2
foreach ($g in $groups) {
    $ws.Cells.Item($r, 2).Value = $displayUrl
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $baseUrl, $gid)
    $ws.Cells.Item($r, 1).Value = $g
    $r++
}

[void]$ws.Range("V6").Select()
